$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was updated
# from 45190 (2023-09-21) to 45192 (2023-09-23) for every data row
# (rows 2 through 135).
for ($r = 2; $r -le 135; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}
